# This script applies a weekly data update to the "Hortaliza, Vega Modelo de
# Temuco - Perejil" sheet: two new observation rows are inserted at the top
# of the data block (just above the former row 173), pushing the existing
# rows 173-187 down to 175-189. The workbook dimension grows from
# A1:R187 to A1:R189 automatically as a result.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before row 173 (existing data shifts down by 2)
$ws.Rows.Item(173).Insert()
$ws.Rows.Item(173).Insert()

# --- New row 173 ---
$ws.Cells.Item(173, 1).Value = 10
$ws.Cells.Item(173, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(173, 3).Value = "La Araucanía"
$ws.Cells.Item(173, 4).Value = 44461
$ws.Cells.Item(173, 4).NumberFormat = $ws.Cells.Item(175, 4).NumberFormat
$ws.Cells.Item(173, 5).Value = 9
$ws.Cells.Item(173, 6).Value = 100112044
$ws.Cells.Item(173, 7).Value = "Perejil"
$ws.Cells.Item(173, 8).Value = "Sin especificar"
$ws.Cells.Item(173, 9).Value = "Primera"
$ws.Cells.Item(173, 10).Value = 30
$ws.Cells.Item(173, 11).Value = 5000
$ws.Cells.Item(173, 12).Value = 5000
$ws.Cells.Item(173, 13).Value = 5000
$ws.Cells.Item(173, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(173, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(173, 16).Value = 1667
$ws.Cells.Item(173, 17).Value = 3
$ws.Cells.Item(173, 18).Value = "Hortaliza"

# --- New row 174 ---
$ws.Cells.Item(174, 1).Value = 10
$ws.Cells.Item(174, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(174, 3).Value = "La Araucanía"
$ws.Cells.Item(174, 4).Value = 44461
$ws.Cells.Item(174, 4).NumberFormat = $ws.Cells.Item(175, 4).NumberFormat
$ws.Cells.Item(174, 5).Value = 9
$ws.Cells.Item(174, 6).Value = 100112044
$ws.Cells.Item(174, 7).Value = "Perejil"
$ws.Cells.Item(174, 8).Value = "Sin especificar"
$ws.Cells.Item(174, 9).Value = "Primera"
$ws.Cells.Item(174, 10).Value = 30
$ws.Cells.Item(174, 11).Value = 3300
$ws.Cells.Item(174, 12).Value = 3300
$ws.Cells.Item(174, 13).Value = 3300
$ws.Cells.Item(174, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(174, 15).Value = "Región Metropolitana"
$ws.Cells.Item(174, 16).Value = 1100
$ws.Cells.Item(174, 17).Value = 3
$ws.Cells.Item(174, 18).Value = "Hortaliza"
